# cn-#18 make sure table rows will not be ignored due to empty space
#
# 1) AddTwoNumbers sheet: move the selection/active cell (no longer the
#    "tabSelected" sheet once SumAverage is added and activated below).
# 2) SumProduct sheet: fix the scenario description text and move the
#    selection to the formula cell.
# 3) Add a brand-new "SumAverage" worksheet (mirrors SumProduct, but
#    demonstrates that empty cells inside the numbers table are still
#    picked up when averaging / summing).

$wb = $excel.ActiveWorkbook

# --- Sheet 1: AddTwoNumbers -------------------------------------------------
$ws1 = $wb.Worksheets.Item("AddTwoNumbers")
$ws1.Range("C5").Select()

# --- Sheet 2: SumProduct ----------------------------------------------------
$ws2 = $wb.Worksheets.Item("SumProduct")
$ws2.Range("B3").Value = "I have following list of numbers"
$ws2.Range("C11").Select()

# --- Sheet 3: SumAverage (new) ----------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws3.Name = "SumAverage"

$ws3.Range("A1").Value = "Scenario"
$ws3.Range("B1").Value = "SumAverage Function"

$ws3.Range("A2").Value = "#"
$ws3.Range("B2").Value = "To test that table value will not be ignored due to empty space"

$ws3.Range("A3").Value = "#"
$ws3.Range("B3").Value = "For simplify testing purpose, I made the roll average full integer."

$ws3.Range("A5").Value = "Given"
$ws3.Range("B5").Value = "I have following list of numbers"

$ws3.Range("B6").Value = "First Number"
$ws3.Range("C6").Value = "Second Number"
$ws3.Range("D6").Value = "Third Number"
$ws3.Range("E6").Value = "Forth Number"
$ws3.Range("H6").Value = "#Avg"

$ws3.Range("B7").Value = 50
$ws3.Range("C7").Value = 70
$ws3.Range("D7").Value = 100
$ws3.Range("E7").Value = 20
$ws3.Range("H7").Formula = "=AVERAGE(B7:E7)"

$ws3.Range("C8").Value = 12
$ws3.Range("D8").Value = 4
$ws3.Range("H8").Formula = "=AVERAGE(B8:E8)"

$ws3.Range("B9").Value = -3
$ws3.Range("E9").Value = 11
$ws3.Range("H9").Formula = "=AVERAGE(B9:E9)"

$ws3.Range("D10").Value = 24
$ws3.Range("E10").Value = 20
$ws3.Range("H10").Formula = "=AVERAGE(B10:E10)"

$ws3.Range("A12").Value = "When"
$ws3.Range("B12").Value = "I call the SumAverage function"

$ws3.Range("A14").Value = "Then"
$ws3.Range("B14").Value = "the result should be"
$ws3.Range("C14").Formula = "=SUM(H7:H10)"

$ws3.Range("A1").EntireColumn.ColumnWidth = 7.666666666666667
$ws3.Range("B1").EntireColumn.ColumnWidth = 28.666666666666668
$ws3.Range("C1").EntireColumn.ColumnWidth = 14.5
$ws3.Range("D1").EntireColumn.ColumnWidth = 12.666666666666666
$ws3.Range("E1").EntireColumn.ColumnWidth = 12.666666666666666

$ws3.Range("C14").Select()
